$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("J5").Value = 0.021739130434782608
$ws.Range("K5").Value = 1

# Row 6
$ws.Range("D6").Value = 0.13157894736842105
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 0.18421052631578946
$ws.Range("G6").Value = 7
$ws.Range("N6").Value = 0.15789473684210525
$ws.Range("O6").Value = 6

# Row 7
$ws.Range("D7").Value = 0.1
$ws.Range("E7").Value = 3

# Row 8
$ws.Range("H8").Value = 0.16071428571428573
$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 0.14285714285714285
$ws.Range("K8").Value = 8
$ws.Range("N8").Value = 0.17857142857142858
$ws.Range("O8").Value = 10

# Row 11
$ws.Range("D11").Value = 0.18181818181818182
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 0.2727272727272727
$ws.Range("G11").Value = 9

# Row 12
$ws.Range("D12").Value = 0.02631578947368421
$ws.Range("E12").Value = 1
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0

# Row 13
$ws.Range("F13").Value = 0.10714285714285714
$ws.Range("G13").Value = 3

# Row 15
$ws.Range("D15").Value = 0.08333333333333333
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.16666666666666666
$ws.Range("G15").Value = 4

# Row 16
$ws.Range("F16").Value = 0.07142857142857142
$ws.Range("G16").Value = 2

# Row 17
$ws.Range("F17").Value = 0.058823529411764705
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 0.35294117647058826
$ws.Range("I17").Value = 6

# Row 18
$ws.Range("D18").Value = 0.2222222222222222
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 0.2777777777777778
$ws.Range("G18").Value = 5

# Row 19
$ws.Range("D19").Value = 0.24390243902439024
$ws.Range("E19").Value = 10

# Row 20
$ws.Range("F20").Value = 0.16326530612244897
$ws.Range("G20").Value = 8
$ws.Range("L20").Value = 0.061224489795918366
$ws.Range("M20").Value = 3
$ws.Range("N20").Value = 0.12244897959183673
$ws.Range("O20").Value = 6

# Row 21
$ws.Range("F21").Value = 0.175
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 0.4
$ws.Range("I21").Value = 16
$ws.Range("N21").Value = 0.125
$ws.Range("O21").Value = 5

# Row 22
$ws.Range("D22").Value = 0.07142857142857142
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 0.11904761904761904
$ws.Range("G22").Value = 5
$ws.Range("L22").Value = 0.11904761904761904
$ws.Range("M22").Value = 5
$ws.Range("N22").Value = 0.21428571428571427
$ws.Range("O22").Value = 9

# Row 26
$ws.Range("D26").Value = 0.0847457627118644
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 0.15254237288135594
$ws.Range("G26").Value = 9
$ws.Range("N26").Value = 0.15254237288135594
$ws.Range("O26").Value = 9

# Row 28
$ws.Range("F28").Value = 0.17391304347826086
$ws.Range("G28").Value = 4

# Row 29
$ws.Range("D29").Value = 0.3
$ws.Range("E29").Value = 6

# Row 31
$ws.Range("N31").Value = 0.15384615384615385
$ws.Range("O31").Value = 2

# Row 32
$ws.Range("H32").Value = 0.08771929824561403
$ws.Range("I32").Value = 5
$ws.Range("L32").Value = 0.07017543859649122
$ws.Range("M32").Value = 4

# Row 34
$ws.Range("N34").Value = 0.030303030303030304
$ws.Range("O34").Value = 1

# Row 36
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("L36").Value = 0.1
$ws.Range("M36").Value = 2
$ws.Range("N36").Value = 0.35
$ws.Range("O36").Value = 7

# Row 38
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0

# Row 40
$ws.Range("J40").Value = 0.06896551724137931
$ws.Range("K40").Value = 2
$ws.Range("N40").Value = 0.10344827586206896
$ws.Range("O40").Value = 3

# Row 41
$ws.Range("D41").Value = 0.05128205128205128
$ws.Range("E41").Value = 2
$ws.Range("H41").Value = 0.1794871794871795
$ws.Range("I41").Value = 7
$ws.Range("L41").Value = 0.1282051282051282
$ws.Range("M41").Value = 5

# Row 42
$ws.Range("F42").Value = 0.058823529411764705
$ws.Range("G42").Value = 4
$ws.Range("L42").Value = 0.058823529411764705
$ws.Range("M42").Value = 4

# Row 43
$ws.Range("F43").Value = 0.04395604395604396
$ws.Range("G43").Value = 4
$ws.Range("H43").Value = 0.07692307692307693
$ws.Range("I43").Value = 7
$ws.Range("N43").Value = 0.0989010989010989
$ws.Range("O43").Value = 9

# Row 44
$ws.Range("H44").Value = 0.05263157894736842
$ws.Range("I44").Value = 1

# Row 45
$ws.Range("F45").Value = 0.3448275862068966
$ws.Range("G45").Value = 10
$ws.Range("H45").Value = 0.4482758620689655
$ws.Range("I45").Value = 13

# Row 46
$ws.Range("J46").Value = 0.09090909090909091
$ws.Range("K46").Value = 2

# Row 47
$ws.Range("D47").Value = 0.05555555555555555
$ws.Range("E47").Value = 2
$ws.Range("H47").Value = 0.1388888888888889
$ws.Range("I47").Value = 5
$ws.Range("N47").Value = 0.08333333333333333
$ws.Range("O47").Value = 3

# Row 48
$ws.Range("H48").Value = 0.10638297872340426
$ws.Range("I48").Value = 5
$ws.Range("L48").Value = 0.14893617021276595
$ws.Range("M48").Value = 7

# Row 49
$ws.Range("J49").Value = 0.018867924528301886
$ws.Range("K49").Value = 1

# Row 53
$ws.Range("F53").Value = 0.027777777777777776
$ws.Range("G53").Value = 1

# Row 54
$ws.Range("H54").Value = 0.13157894736842105
$ws.Range("I54").Value = 5
$ws.Range("N54").Value = 0.10526315789473684
$ws.Range("O54").Value = 4

# Row 55
$ws.Range("F55").Value = 0.2413793103448276
$ws.Range("G55").Value = 7
$ws.Range("H55").Value = 0.2413793103448276
$ws.Range("I55").Value = 7
$ws.Range("N55").Value = 0.06896551724137931
$ws.Range("O55").Value = 2

